$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.27738566666667
$ws.Range("H2").Value = 30.832157
$ws.Range("I2").Value = 0.3571200664977529
$ws.Range("J2").Value = 0.3571200664977529
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 24.55401366666667
$ws.Range("N2").Value = 73.662041
$ws.Range("O2").Value = 0.2501033965205259
$ws.Range("P2").Value = 0.250103396520526
$ws.Range("Q2").Value = 252.3510681169374
$ws.Range("R2").Value = 2271.159613052437
$ws.Range("S2").Value = 0.08931694159672407
$ws.Range("T2").Value = 0.0893169415967241

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.27738566666667
$ws.Range("H3").Value = 30.832157
$ws.Range("I3").Value = 0.3571200664977529
$ws.Range("J3").Value = 0.3571200664977529
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.205204666666667
$ws.Range("N3").Value = 3.615614
$ws.Range("O3").Value = 0.01227602886956614
$ws.Range("P3").Value = 0.01227602886956614
$ws.Range("Q3").Value = 12.38635316659978
$ws.Range("R3").Value = 111.477178499398
$ws.Range("S3").Value = 0.004384016246227794
$ws.Range("T3").Value = 0.004384016246227796

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.27738566666667
$ws.Range("H4").Value = 30.832157
$ws.Range("I4").Value = 0.3571200664977529
$ws.Range("J4").Value = 0.3571200664977529
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 67.337282
$ws.Range("N4").Value = 202.011846
$ws.Range("O4").Value = 0.6858871697837075
$ws.Range("P4").Value = 0.6858871697837076
$ws.Range("Q4").Value = 692.0512168590913
$ws.Range("R4").Value = 6228.460951731821
$ws.Range("S4").Value = 0.2449440716831131
$ws.Range("T4").Value = 0.2449440716831132

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 10.27738566666667
$ws.Range("H5").Value = 30.832157
$ws.Range("I5").Value = 0.3571200664977529
$ws.Range("J5").Value = 0.3571200664977529
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.078950333333333
$ws.Range("N5").Value = 15.236851
$ws.Range("O5").Value = 0.0517334048262004
$ws.Range("P5").Value = 0.05173340482620041
$ws.Range("Q5").Value = 52.19833135751188
$ws.Range("R5").Value = 469.784982217607
$ws.Range("S5").Value = 0.01847503697168786
$ws.Range("T5").Value = 0.01847503697168786

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.987132666666667
$ws.Range("H6").Value = 5.961398
$ws.Range("I6").Value = 0.06904917000064482
$ws.Range("J6").Value = 0.06904917000064482
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 24.55401366666667
$ws.Range("N6").Value = 73.662041
$ws.Range("O6").Value = 0.2501033965205259
$ws.Range("P6").Value = 0.250103396520526
$ws.Range("Q6").Value = 48.79208265481311
$ws.Range("R6").Value = 439.128743893318
$ws.Range("S6").Value = 0.01726943194408448
$ws.Range("T6").Value = 0.01726943194408448

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.987132666666667
$ws.Range("H7").Value = 5.961398
$ws.Range("I7").Value = 0.06904917000064482
$ws.Range("J7").Value = 0.06904917000064482
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.205204666666667
$ws.Range("N7").Value = 3.615614
$ws.Range("O7").Value = 0.01227602886956614
$ws.Range("P7").Value = 0.01227602886956614
$ws.Range("Q7").Value = 2.394901563152445
$ws.Range("R7").Value = 21.554114068372
$ws.Range("S7").Value = 0.0008476496043474963
$ws.Range("T7").Value = 0.0008476496043474964

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.987132666666667
$ws.Range("H8").Value = 5.961398
$ws.Range("I8").Value = 0.06904917000064482
$ws.Range("J8").Value = 0.06904917000064482
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 67.337282
$ws.Range("N8").Value = 202.011846
$ws.Range("O8").Value = 0.6858871697837075
$ws.Range("P8").Value = 0.6858871697837076
$ws.Range("Q8").Value = 133.8081127467453
$ws.Range("R8").Value = 1204.273014720708
$ws.Range("S8").Value = 0.04735993978765635
$ws.Range("T8").Value = 0.04735993978765636

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.987132666666667
$ws.Range("H9").Value = 5.961398
$ws.Range("I9").Value = 0.06904917000064482
$ws.Range("J9").Value = 0.06904917000064482
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.078950333333333
$ws.Range("N9").Value = 15.236851
$ws.Range("O9").Value = 0.0517334048262004
$ws.Range("P9").Value = 0.05173340482620041
$ws.Range("Q9").Value = 10.09254811974422
$ws.Range("R9").Value = 90.832933077698
$ws.Range("S9").Value = 0.003572148664556491
$ws.Range("T9").Value = 0.003572148664556492

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 15.70994033333333
$ws.Range("H10").Value = 47.129821
$ws.Range("I10").Value = 0.5458912527445677
$ws.Range("J10").Value = 0.5458912527445677
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 24.55401366666667
$ws.Range("N10").Value = 73.662041
$ws.Range("O10").Value = 0.2501033965205259
$ws.Range("P10").Value = 0.250103396520526
$ws.Range("Q10").Value = 385.7420896471846
$ws.Range("R10").Value = 3471.678806824661
$ws.Range("S10").Value = 0.1365292564422613
$ws.Range("T10").Value = 0.1365292564422613

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 15.70994033333333
$ws.Range("H11").Value = 47.129821
$ws.Range("I11").Value = 0.5458912527445677
$ws.Range("J11").Value = 0.5458912527445677
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.205204666666667
$ws.Range("N11").Value = 3.615614
$ws.Range("O11").Value = 0.01227602886956614
$ws.Range("P11").Value = 0.01227602886956614
$ws.Range("Q11").Value = 18.93369340278822
$ws.Range("R11").Value = 170.403240625094
$ws.Range("S11").Value = 0.006701376778335941
$ws.Range("T11").Value = 0.006701376778335942

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 15.70994033333333
$ws.Range("H12").Value = 47.129821
$ws.Range("I12").Value = 0.5458912527445677
$ws.Range("J12").Value = 0.5458912527445677
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 67.337282
$ws.Range("N12").Value = 202.011846
$ws.Range("O12").Value = 0.6858871697837075
$ws.Range("P12").Value = 0.6858871697837076
$ws.Range("Q12").Value = 1057.864682428841
$ws.Range("R12").Value = 9520.782141859565
$ws.Range("S12").Value = 0.3744198063546541
$ws.Range("T12").Value = 0.3744198063546542

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 15.70994033333333
$ws.Range("H13").Value = 47.129821
$ws.Range("I13").Value = 0.5458912527445677
$ws.Range("J13").Value = 0.5458912527445677
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 5.078950333333333
$ws.Range("N13").Value = 15.236851
$ws.Range("O13").Value = 0.0517334048262004
$ws.Range("P13").Value = 0.05173340482620041
$ws.Range("Q13").Value = 79.79000669263012
$ws.Range("R13").Value = 718.110060233671
$ws.Range("S13").Value = 0.0282408131693164
$ws.Range("T13").Value = 0.02824081316931641

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.8040576666666667
$ws.Range("H14").Value = 2.412173
$ws.Range("I14").Value = 0.02793951075703474
$ws.Range("J14").Value = 0.02793951075703475
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 24.55401366666667
$ws.Range("N14").Value = 73.662041
$ws.Range("O14").Value = 0.2501033965205259
$ws.Range("P14").Value = 0.250103396520526
$ws.Range("Q14").Value = 19.74284293612144
$ws.Range("R14").Value = 177.685586425093
$ws.Range("S14").Value = 0.006987766537456161
$ws.Range("T14").Value = 0.006987766537456163

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.8040576666666667
$ws.Range("H15").Value = 2.412173
$ws.Range("I15").Value = 0.02793951075703474
$ws.Range("J15").Value = 0.02793951075703475
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.205204666666667
$ws.Range("N15").Value = 3.615614
$ws.Range("O15").Value = 0.01227602886956614
$ws.Range("P15").Value = 0.01227602886956614
$ws.Range("Q15").Value = 0.9690540521357778
$ws.Range("R15").Value = 8.721486469222
$ws.Range("S15").Value = 0.0003429862406549123
$ws.Range("T15").Value = 0.0003429862406549124

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.8040576666666667
$ws.Range("H16").Value = 2.412173
$ws.Range("I16").Value = 0.02793951075703474
$ws.Range("J16").Value = 0.02793951075703475
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 67.337282
$ws.Range("N16").Value = 202.011846
$ws.Range("O16").Value = 0.6858871697837075
$ws.Range("P16").Value = 0.6858871697837076
$ws.Range("Q16").Value = 54.14305784459533
$ws.Range("R16").Value = 487.287520601358
$ws.Range("S16").Value = 0.01916335195828401
$ws.Range("T16").Value = 0.01916335195828402

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.8040576666666667
$ws.Range("H17").Value = 2.412173
$ws.Range("I17").Value = 0.02793951075703474
$ws.Range("J17").Value = 0.02793951075703475
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 5.078950333333333
$ws.Range("N17").Value = 15.236851
$ws.Range("O17").Value = 0.0517334048262004
$ws.Range("P17").Value = 0.05173340482620041
$ws.Range("Q17").Value = 4.083768954135889
$ws.Range("R17").Value = 36.753920587223
$ws.Range("S17").Value = 0.001445406020639659
$ws.Range("T17").Value = 0.00144540602063966
